$wb = $excel.ActiveWorkbook

# ----- Sheet "Desired Data" (first sheet) -----
$ws1 = $wb.Worksheets.Item("Desired Data")

# Row 1: label cell R1 changes text
$ws1.Range("R1").Value = "interval11_count"

# Row 2: header row for interval columns.
# Shift labels so P starts at interval1_count ... AA ends at interval12_count,
# dropping interval0_count and adding interval12_count. Columns S:Z keep their
# existing (unshifted) cell contents untouched, matching the original diff.
$ws1.Range("P2").Value = "interval1_count"
$ws1.Range("P2").Interior.Pattern = -4142   # drop fill so style matches Q2:AA2 (no-fill bordered header)
$ws1.Range("Q2").Value = "interval2_count"
$ws1.Range("R2").Value = "interval3_count"
$ws1.Range("AA2").Value = "interval12_count"
$ws1.Range("AB2").Clear()

# Row 3: placeholder values become "shape = (1, )" with text-format style (no border flag)
$ws1.Range("P3:AA3").Value = "shape = (1, )"
$ws1.Range("P3:AA3").NumberFormat = "@"
$ws1.Range("AB3").Clear()

# Row 4: same placeholder text, default (no) style
$ws1.Range("P4:AA4").Value = "shape = (1, )"
$ws1.Range("P4:AA4").ClearFormats()
$ws1.Range("AB4").Clear()

# Row 5: same as row 4
$ws1.Range("P5:AA5").Value = "shape = (1, )"
$ws1.Range("P5:AA5").ClearFormats()
$ws1.Range("AB5").Clear()

# Row 6: entire block removed
$ws1.Range("P6:AB6").Clear()

# Row 7: only P7 keeps a (plain-styled) value, rest removed
$ws1.Range("P7").Value = "Process for calculating interval"
$ws1.Range("P7").ClearFormats()
$ws1.Range("Q7:AB7").Clear()

# Row 8
$ws1.Range("P8").Value = "1. Find highest probability note for each row. Result: column of size = (971, ) containing integer of the index location of highest probability note"
$ws1.Range("P8").ClearFormats()
$ws1.Range("Q8:AB8").Clear()

# Row 9
$ws1.Range("P9").Value = "2. Use the index locations to determine intervals"
$ws1.Range("P9").ClearFormats()
$ws1.Range("Q9:AB9").Clear()

# Row 10
$ws1.Range("P10").Value = "2.1 interval[n] = abs(index [n+1] - index[n])"
$ws1.Range("P10").ClearFormats()
$ws1.Range("Q10:AB10").Clear()

# Row 11
$ws1.Range("P11").Value = "3. Parse interval array into columns based on how many times each interval appears"
$ws1.Range("P11").ClearFormats()
$ws1.Range("Q11:AB11").Clear()

# Row 12: entire block removed
$ws1.Range("P12:AB12").Clear()

# Row 13: entire block removed
$ws1.Range("P13:AB13").Clear()

# ----- Sheet "Sheet1" (second sheet) -----
$ws2 = $wb.Worksheets.Item("Sheet1")
$ws2.Range("C15").Value = "Album/Single/EP"
$ws2.Range("D15").Value = "Song title"

# ----- View / selection state -----
$ws1.Range("P1:P10").Select()
$ws1.Range("P1").Activate()
$ws1.Application.ActiveWindow.ScrollRow = 1
$ws1.Application.ActiveWindow.ScrollColumn = 2

$ws2.Range("B16:H26").Select()
$ws2.Range("B16").Activate()
$ws2.Application.ActiveWindow.Zoom = 130

$ws2.Activate()
